$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the corrected LED part (replacing the old blue LED entry
# with the new 2-PLCC Orange LED entry)
$ws.Range("A2").Value = "LO T67K-K1L2-24-Z"
$ws.Range("B2").Value = "475-2745-1-ND"
$ws.Range("D2").Value = "2-PLCC Orange LED"

# Update the active selection to D2, matching the saved view state
$ws.Range("D2").Select()
